# Insert a new row for "Citation to discover assets to debtor" into the
# alphabetically-sorted list on Sheet1, immediately after "Cannabis
# expungement" (row 4) and before "Civil No Contact Order - CNCO"
# (previously row 5), then restore hyperlinks on the rows that moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Insert a new row at position 5 - this shifts existing rows 5:48 down
#    to 6:49 and carries formatting (incl. the Hyperlink style on column B)
#    down with them.
$ws.Rows.Item(5).Insert()

# 2. Populate the new row.
$ws.Cells.Item(5, 1).Value = "Citation to discover assets to debtor"
$ws.Cells.Item(5, 2).Value = "https://www.illinoislegalaid.org/legal-information/citation-discover-assets-debtor"

# 3. The row insert does not carry the worksheet's Hyperlinks collection
#    along with it, so rebuild it: drop the stale entries and re-add a
#    hyperlink for every row (now shifted down by one, where applicable)
#    that had one before the insert. The new row 5 intentionally gets no
#    hyperlink, matching the rest of the freshly-added rows in this sheet
#    that only have the visual Hyperlink style without a live link.
$ws.Hyperlinks.Delete()

$linkedRows = @(2, 4, 7, 8, 11, 12, 13, 14, 15, 18, 19, 20, 21, 23, 24, 25, 27, 30, 31, 32, 33, 35, 36, 37, 38, 40, 41, 42, 43, 44, 46, 47, 48, 49)
foreach ($r in $linkedRows) {
    $cell = $ws.Cells.Item($r, 2)
    $ws.Hyperlinks.Add($cell, $cell.Value())
    # Hyperlinks.Add silently reassigns a fresh cell style; put the
    # worksheet's existing "Hyperlink" style back so column B keeps using
    # the same style index (s="1") it used before this script ran.
    $cell.Style = "Hyperlink"
}

# 4. The sheet carries a stale sortState/sortCondition (left over from an
#    earlier manual sort of the first block of rows) that only covered
#    A2:B17 before the insert; re-apply the same sort over its new extent
#    (A2:B18) so that bookkeeping grows by one row along with everything
#    else that shifted.
$sortRange = $ws.Range("A2:B18")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A18")) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# 5. Mirror the author's final selection.
$ws.Range("A50").Select()
